$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.223.35"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.174.10"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.07%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.85"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.172.24"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.99%  "

$ws.Range("E9").Value = "  +1.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.67"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -7.70%  "

$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("E13").Value = "  -2.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.64"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.699.13"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.340.25"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.38"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.177.78"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("E21").Value = "  -2.15%  "

$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("E23").Value = "  +1.91%  "

$ws.Range("E24").Value = "  -2.15%  "

$ws.Range("E25").Value = "  -1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.15"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.14%  "

$ws.Range("E29").Value = "  +6.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.07"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("E34").Value = "  -1.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.52"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "505.93"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.66"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0892"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.39%  "

$ws.Range("E39").Value = "  -0.82%  "

$ws.Range("E40").Value = "  +5.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.85"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.299"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0671"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.840.44"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.22"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.66%  "

$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.54"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.10%  "
